$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 937.5333000000001
$ws.Range("I19").Value = 718.5
$ws.Range("K19").Value = 718.5
$ws.Range("M19").Value = -543.5

$ws.Range("H55").Value = 417.25
$ws.Range("I55").Value = 389.66666
$ws.Range("K55").Value = 389.66666
$ws.Range("M55").Value = -175.66666

$ws.Range("H69").Value = 5222.727
$ws.Range("J69").Value = 4644.2856
$ws.Range("L69").Value = 13932.8568
$ws.Range("N69").Value = -15680.8568

$ws.Range("H72").Value = 5222.727
$ws.Range("J72").Value = 4644.2856
$ws.Range("L72").Value = 41798.5704
$ws.Range("N72").Value = -50534.5704

$ws.Range("H125").Value = 1134.625
$ws.Range("I125").Value = 846.1667
$ws.Range("J125").Value = 2000
$ws.Range("K125").Value = 7615.5003
$ws.Range("L125").Value = 18000
$ws.Range("M125").Value = -5155.5003
$ws.Range("N125").Value = -22920

$ws.Range("H132").Value = 3985.705
$ws.Range("I132").Value = 3074.7585
$ws.Range("J132").Value = 6627.45
$ws.Range("K132").Value = 9224.2755
$ws.Range("L132").Value = 19882.35
$ws.Range("M132").Value = -6694.2755
$ws.Range("N132").Value = -24942.35

$ws.Range("H138").Value = 2242.4075
$ws.Range("I138").Value = 1271.8379
$ws.Range("J138").Value = 4354.8237
$ws.Range("K138").Value = 3815.5137
$ws.Range("L138").Value = 13064.4711
$ws.Range("M138").Value = 1324.4863
$ws.Range("N138").Value = -23344.4711

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 2001.2
$ws.Range("I33").Value = 1244.25
$ws.Range("J33").Value = 5029
$ws.Range("K33").Value = 1244.25
$ws.Range("L33").Value = 5029
$ws.Range("M33").Value = -915.25
$ws.Range("N33").Value = -5687

$ws.Range("H61").Value = 2126.75
$ws.Range("I61").Value = 1355.625
$ws.Range("J61").Value = 2787.7144
$ws.Range("K61").Value = 1355.625
$ws.Range("L61").Value = 2787.7144
$ws.Range("M61").Value = -1143.625
$ws.Range("N61").Value = -3211.7144

$ws.Range("H132").Value = 4020.7942
$ws.Range("I132").Value = 2762.6924
$ws.Range("J132").Value = 5712.724
$ws.Range("K132").Value = 8288.0772
$ws.Range("L132").Value = 17138.172
$ws.Range("M132").Value = -5758.0772
$ws.Range("N132").Value = -22198.172

$ws.Range("H136").Value = 2126.75
$ws.Range("I136").Value = 1355.625
$ws.Range("J136").Value = 2787.7144
$ws.Range("K136").Value = 4066.875
$ws.Range("L136").Value = 8363.143199999999
$ws.Range("M136").Value = -1516.875
$ws.Range("N136").Value = -13463.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2612.2222
$ws.Range("I105").Value = 1522
$ws.Range("J105").Value = 3975
$ws.Range("K105").Value = 1522
$ws.Range("L105").Value = 3975
$ws.Range("M105").Value = 225
$ws.Range("N105").Value = -7469

$ws.Range("H107").Value = 3012.2
$ws.Range("I107").Value = 3015.25
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 3015.25
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -1095.25
$ws.Range("N107").Value = -6840

$ws.Range("H134").Value = 5922.3413
$ws.Range("I134").Value = 3019.25
$ws.Range("J134").Value = 7780.32
$ws.Range("K134").Value = 9057.75
$ws.Range("L134").Value = 23340.96
$ws.Range("M134").Value = -6522.75
$ws.Range("N134").Value = -28410.96

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 980.9
$ws.Range("I107").Value = 312.91666
$ws.Range("J107").Value = 1982.875
$ws.Range("K107").Value = 312.91666
$ws.Range("L107").Value = 1982.875
$ws.Range("M107").Value = 1607.08334
$ws.Range("N107").Value = -5822.875

$ws.Range("H122").Value = 3880
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3880
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 11640
$ws.Range("N122").Value = -16540
$ws.Range("M122").ClearContents()

$ws.Range("H134").Value = 1600.2106
$ws.Range("I134").Value = 939.36365
$ws.Range("J134").Value = 1869.4445
$ws.Range("K134").Value = 2818.09095
$ws.Range("L134").Value = 5608.333500000001
$ws.Range("M134").Value = -283.0909499999998
$ws.Range("N134").Value = -10678.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 839
$ws.Range("I68").Value = 375
$ws.Range("K68").Value = 1125
$ws.Range("M68").Value = -314

$ws.Range("H71").Value = 839
$ws.Range("I71").Value = 375
$ws.Range("K71").Value = 3375
$ws.Range("M71").Value = 681

$ws.Range("H132").Value = 2408
$ws.Range("I132").Value = 937.25
$ws.Range("J132").Value = 3878.75
$ws.Range("K132").Value = 8435.25
$ws.Range("L132").Value = 34908.75
$ws.Range("M132").Value = -5905.25
$ws.Range("N132").Value = -39968.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 77.75
$ws.Range("I2").Value = 70.333336
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 70.333336
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = 42.666664
$ws.Range("N2").Value = -326

$ws.Range("H43").Value = 1000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1839.3793
$ws.Range("I82").Value = 1262.2106
$ws.Range("J82").Value = 2936
$ws.Range("K82").Value = 1262.2106
$ws.Range("L82").Value = 2936
$ws.Range("M82").Value = -901.2106000000001
$ws.Range("N82").Value = -3658

$ws.Range("H85").Value = 1839.3793
$ws.Range("I85").Value = 1262.2106
$ws.Range("J85").Value = 2936
$ws.Range("K85").Value = 1262.2106
$ws.Range("L85").Value = 2936
$ws.Range("M85").Value = -14.21060000000011
$ws.Range("N85").Value = -5432

$ws.Range("H132").Value = 32261340
$ws.Range("I132").Value = 43481188
$ws.Range("J132").Value = 4276
$ws.Range("K132").Value = 130443564
$ws.Range("L132").Value = 12828
$ws.Range("M132").Value = -130441034
$ws.Range("N132").Value = -17888

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2120.3953
$ws.Range("I126").Value = 2008.6129
$ws.Range("J126").Value = 2409.1667
$ws.Range("K126").Value = 6025.8387
$ws.Range("L126").Value = 7227.500100000001
$ws.Range("M126").Value = -3555.8387
$ws.Range("N126").Value = -12167.5001
